$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G20").Value = 39.787234040000001
$ws.Range("H20").Value = 132.49148940000001
$ws.Range("I20").Value = 264.98297869999999
$ws.Range("J20").Value = 405.6382979
$ws.Range("K20").Value = 554.04255320000004
$ws.Range("L20").Value = 702.44680849999997
$ws.Range("M20").Value = 851.06382980000001
$ws.Range("AA21").Value = 851.06382980000001
$ws.Range("G21").Value = 17.937649879999999
$ws.Range("H21").Value = 59.732374100000001
$ws.Range("I21").Value = 119.4647482
$ws.Range("J21").Value = 182.87769779999999
$ws.Range("K21").Value = 249.7841727
$ws.Range("L21").Value = 316.69064750000001
$ws.Range("M21").Value = 383.6930456
$ws.Range("U21").Value = 39.787234040000001
$ws.Range("V21").Value = 132.49148940000001
$ws.Range("W21").Value = 264.98297869999999
$ws.Range("X21").Value = 405.6382979
$ws.Range("Y21").Value = 554.04255320000004
$ws.Range("Z21").Value = 702.44680849999997
$ws.Range("AA22").Value = 383.6930456
$ws.Range("U22").Value = 17.937649879999999
$ws.Range("V22").Value = 59.732374100000001
$ws.Range("W22").Value = 119.4647482
$ws.Range("X22").Value = 182.87769779999999
$ws.Range("Y22").Value = 249.7841727
$ws.Range("Z22").Value = 316.69064750000001
$ws.Range("G25").Value = 10.56497175
$ws.Range("H25").Value = 35.181355930000002
$ws.Range("I25").Value = 70.362711860000005
$ws.Range("J25").Value = 107.7118644
$ws.Range("K25").Value = 147.11864410000001
$ws.Range("L25").Value = 186.5254237
$ws.Range("M25").Value = 225.98870059999999
$ws.Range("AA26").Value = 225.98870059999999
$ws.Range("G26").Value = 5.9935897440000003
$ws.Range("H26").Value = 19.958653850000001
$ws.Range("I26").Value = 39.917307690000001
$ws.Range("J26").Value = 61.10576923
$ws.Range("K26").Value = 83.46153846
$ws.Range("L26").Value = 105.8173077
$ws.Range("M26").Value = 128.20512819999999
$ws.Range("U26").Value = 10.56497175
$ws.Range("V26").Value = 35.181355930000002
$ws.Range("W26").Value = 70.362711860000005
$ws.Range("X26").Value = 107.7118644
$ws.Range("Y26").Value = 147.11864410000001
$ws.Range("Z26").Value = 186.5254237
$ws.Range("AA27").Value = 128.20512819999999
$ws.Range("U27").Value = 5.9935897440000003
$ws.Range("V27").Value = 19.958653850000001
$ws.Range("W27").Value = 39.917307690000001
$ws.Range("X27").Value = 61.10576923
$ws.Range("Y27").Value = 83.46153846
$ws.Range("Z27").Value = 105.8173077
$ws.Range("G31").Value = 25.704467350000002
$ws.Range("H31").Value = 85.595876290000007
$ws.Range("I31").Value = 171.1917526
$ws.Range("J31").Value = 262.06185570000002
$ws.Range("K31").Value = 357.93814429999998
$ws.Range("L31").Value = 453.81443300000001
$ws.Range("M31").Value = 549.82817869999997
$ws.Range("AA32").Value = 549.82817869999997
$ws.Range("G32").Value = 42.25988701
$ws.Range("H32").Value = 140.72542369999999
$ws.Range("I32").Value = 281.45084750000001
$ws.Range("J32").Value = 430.84745759999998
$ws.Range("K32").Value = 588.47457629999997
$ws.Range("L32").Value = 746.10169489999998
$ws.Range("M32").Value = 903.95480229999998
$ws.Range("U32").Value = 25.704467350000002
$ws.Range("V32").Value = 85.595876290000007
$ws.Range("W32").Value = 171.1917526
$ws.Range("X32").Value = 262.06185570000002
$ws.Range("Y32").Value = 357.93814429999998
$ws.Range("Z32").Value = 453.81443300000001
$ws.Range("AA33").Value = 903.95480229999998
$ws.Range("G33").Value = 8.3111111110000007
$ws.Range("H33").Value = 27.675999999999998
$ws.Range("I33").Value = 55.351999999999997
$ws.Range("J33").Value = 84.733333329999994
$ws.Range("K33").Value = 115.7333333
$ws.Range("L33").Value = 146.7333333
$ws.Range("M33").Value = 177.7777778
$ws.Range("U33").Value = 42.25988701
$ws.Range("V33").Value = 140.72542369999999
$ws.Range("W33").Value = 281.45084750000001
$ws.Range("X33").Value = 430.84745759999998
$ws.Range("Y33").Value = 588.47457629999997
$ws.Range("Z33").Value = 746.10169489999998
$ws.Range("AA34").Value = 177.7777778
$ws.Range("G34").Value = 6.1513157889999999
$ws.Range("H34").Value = 20.483881579999998
$ws.Range("I34").Value = 40.967763159999997
$ws.Range("J34").Value = 62.713815789999998
$ws.Range("K34").Value = 85.657894740000003
$ws.Range("L34").Value = 108.6019737
$ws.Range("M34").Value = 131.5789474
$ws.Range("U34").Value = 8.3111111110000007
$ws.Range("V34").Value = 27.675999999999998
$ws.Range("W34").Value = 55.351999999999997
$ws.Range("X34").Value = 84.733333329999994
$ws.Range("Y34").Value = 115.7333333
$ws.Range("Z34").Value = 146.7333333
$ws.Range("AA35").Value = 131.5789474
$ws.Range("U35").Value = 6.1513157889999999
$ws.Range("V35").Value = 20.483881579999998
$ws.Range("W35").Value = 40.967763159999997
$ws.Range("X35").Value = 62.713815789999998
$ws.Range("Y35").Value = 85.657894740000003
$ws.Range("Z35").Value = 108.6019737
$ws.Range("G41").Value = 119.68
$ws.Range("H41").Value = 398.53440000000001
$ws.Range("I41").Value = 797.06880000000001
$ws.Range("J41").Value = 1220.1600000000001
$ws.Range("K41").Value = 1666.56
$ws.Range("L41").Value = 2112.96
$ws.Range("M41").Value = 2560
$ws.Range("AA42").Value = 2560
$ws.Range("U42").Value = 119.68
$ws.Range("V42").Value = 398.53440000000001
$ws.Range("W42").Value = 797.06880000000001
$ws.Range("X42").Value = 1220.1600000000001
$ws.Range("Y42").Value = 1666.56
$ws.Range("Z42").Value = 2112.96

$ws.Range("S27").Select() | Out-Null
